$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.054961085319519
$ws.Range("B1").Value = 2.41200065612793
$ws.Range("C1").Value = 5.152599334716797
$ws.Range("D1").Value = 2.314255237579346
$ws.Range("E1").Value = 1.315142869949341
